$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the Neo4j query text stored in cell B4 (FilesTab query):
# Remove the "File Type" coalesce line and the "Breed" coalesce line,
# matching the edit recorded in the shared string table.
$newQuery = "MATCH (f:file)-->(parent)`n" +
    "WITH DISTINCT f, parent`n" +
    "MATCH (f)-[*]->(c:case)<--(demo:demographic)`n" +
    " MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`n" +
    "WHERE s.clinical_study_designation IN ['COTC007B']`n" +
    "WITH DISTINCT f, parent, c, demo, diag, s`n" +
    "RETURN coalesce(f.file_name, '') AS ``File Name``, `n" +
    "        coalesce(labels(parent)[0], '') AS ``Association``,`n" +
    "        coalesce(f.file_description, '') AS ``Description``,`n" +
    "        coalesce(f.file_format, '') AS ``Format``,`n" +
    "        coalesce(f.file_size, '') AS ``Size``,`n" +
    "        coalesce(c.case_id, '') AS ``Case ID``, `n" +
    "        coalesce(diag.disease_term,'') AS Diagnosis , `n" +
    "        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$ws.Range("B4").Value = $newQuery

# The row shrank now that it holds one less line of wrapped text.
$ws.Rows.Item(4).RowHeight = 203

# The sheet had scrolled down one row and the selection moved from B3 to B4.
$ws.Range("B4").Select()
